$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'247.20"
$ws.Range("D3").Value = "'26.37"
$ws.Range("D4").Value = "'5.084"
$ws.Range("D5").Value = "'0.05621"
$ws.Range("D6").Value = "'6.520"
$ws.Range("D7").Value = "'0.8134"
$ws.Range("D8").Value = "'0.8492"
$ws.Range("B9").Value = "'One"
$ws.Range("C9").Value = "'https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D9").Value = "'0.009876"
$ws.Range("E9").Value = "'8OneONEBestin24h"
$ws.Range("B10").Value = "'LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "'0.03207"
$ws.Range("E10").Value = "'9LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("B11").Value = "'BitrueCoin"
$ws.Range("C11").Value = "'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D11").Value = "'0.02821"
$ws.Range("E11").Value = "'10BitrueCoinBTR"
$ws.Range("B12").Value = "'BitMartToken"
$ws.Range("C12").Value = "'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D12").Value = "'0.09408"
$ws.Range("E12").Value = "'11BitMartTokenBMX"
$ws.Range("B13").Value = "'BitForexToken"
$ws.Range("C13").Value = "'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D13").Value = "'0.001511"
$ws.Range("E13").Value = "'12BitForexTokenBF"
$ws.Range("B14").Value = "'TigerCash"
$ws.Range("C14").Value = "'https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D14").Value = "'0.006128"
$ws.Range("E14").Value = "'13TigerCashTCH"
$ws.Range("B15").Value = "'LEO"
$ws.Range("C15").Value = "'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D15").Value = "'3.587"
$ws.Range("E15").Value = "'14LEOLEO"
$ws.Range("B16").Value = "'GateToken"
$ws.Range("C16").Value = "'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D16").Value = "'3.060"
$ws.Range("E16").Value = "'15GateTokenGT"
$ws.Range("B17").Value = "'BTSEToken"
$ws.Range("C17").Value = "'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D17").Value = "'2.118"
$ws.Range("E17").Value = "'16BTSETokenBTSE"
$ws.Range("B18").Value = "'BitpandaEcosystemToken"
$ws.Range("C18").Value = "'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D18").Value = "'0.3181"
$ws.Range("E18").Value = "'17BitpandaEcosystemTokenBEST"
$ws.Range("B19").Value = "'WazirX"
$ws.Range("C19").Value = "'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D19").Value = "'0.1348"
$ws.Range("E19").Value = "'18WazirXWRX"
$ws.Range("B20").Value = "'MandalaExchangeToken"
$ws.Range("C20").Value = "'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D20").Value = "'0.06963"
$ws.Range("E20").Value = "'19MandalaExchangeTokenMDX"
$ws.Range("B21").Value = "'ProBitToken"
$ws.Range("C21").Value = "'https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D21").Value = "'0.1320"
$ws.Range("E21").Value = "'20ProBitTokenPROB"
$ws.Range("B22").Value = "'MCDex"
$ws.Range("C22").Value = "'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D22").Value = "'3.739"
$ws.Range("E22").Value = "'21MCDexMCB"
$ws.Range("B23").Value = "'CoinExToken"
$ws.Range("C23").Value = "'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D23").Value = "'0.04647"
$ws.Range("E23").Value = "'22CoinExTokenCET"
$ws.Range("B24").Value = "'ZBToken"
$ws.Range("C24").Value = "'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D24").Value = "'0.1350"
$ws.Range("E24").Value = "'23ZBTokenZB"
$ws.Range("D25").Value = "'0.001248"
$ws.Range("D26").Value = "'0.004612"
$ws.Range("B41").Value = "'KickToken"
$ws.Range("C41").Value = "'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").Value = "'0.006150"
$ws.Range("E41").Value = "'40KickTokenKICK"
$ws.Range("B42").Value = "'BKEXToken"
$ws.Range("C42").Value = "'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "'0.1060"
$ws.Range("E42").Value = "'41BKEXTokenBKK"
$ws.Range("D43").Value = "'0.002500"
$ws.Range("D44").Value = "'0.008665"
$ws.Range("D45").Value = "'0.00005294"
$ws.Range("D48").Value = "'0.002555"
